$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the amount of Chillers (D2, under "Chillers" / "Amount:")
$ws.Range("D2").Value = 4

# Fix the room coordinate for the fourth corner (A5): was (4,0), now (2,0)
$ws.Range("A5").Value = "(2,0)"

# Update the selected cell to reflect where the user was working
$ws.Range("D3").Select()
